$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.812326845052874
$ws.Range("D2").Value = 8.153630033653448
$ws.Range("E2").Value = 13.25215081241025
$ws.Range("F2").Value = 37.72950447095415
$ws.Range("G2").Value = 43.07830343207233
$ws.Range("H2").Value = 17.66256062012796
$ws.Range("J2").Value = 10.0955432465944
$ws.Range("K2").Value = 14.73382815716467
$ws.Range("N2").Value = 19.46770963653784
$ws.Range("B3").Value = 7.741325745014461
$ws.Range("D3").Value = 8.121171058301368
$ws.Range("E3").Value = 13.2045482296512
$ws.Range("F3").Value = 37.68549111550069
$ws.Range("G3").Value = 42.95113738155442
$ws.Range("H3").Value = 17.696518927952
$ws.Range("J3").Value = 10.10302830563148
$ws.Range("K3").Value = 14.40340377692292
$ws.Range("N3").Value = 19.53797876002085
$ws.Range("B4").Value = 7.699263296378494
$ws.Range("D4").Value = 8.102459414656478
$ws.Range("E4").Value = 13.17805137071585
$ws.Range("F4").Value = 37.66876990377961
$ws.Range("G4").Value = 42.88717265147123
$ws.Range("H4").Value = 17.72103222835331
$ws.Range("J4").Value = 10.10934057483708
$ws.Range("K4").Value = 14.19989791276226
$ws.Range("N4").Value = 19.58307156400298
$ws.Range("B5").Value = 7.682526771676789
$ws.Range("D5").Value = 8.09514587106082
$ws.Range("E5").Value = 13.1679479600762
$ws.Range("F5").Value = 37.66454846781727
$ws.Range("G5").Value = 42.86466852709383
$ws.Range("H5").Value = 17.73194042917924
$ws.Range("J5").Value = 10.11234443379036
$ws.Range("K5").Value = 14.11693451764509
$ws.Range("N5").Value = 19.60193855896523
$ws.Range("B6").Value = 7.679772638576911
$ws.Range("D6").Value = 8.093950429609306
$ws.Range("E6").Value = 13.16631243118732
$ws.Range("F6").Value = 37.66400409513754
$ws.Range("G6").Value = 42.86114714188293
$ws.Range("H6").Value = 17.73380715491813
$ws.Range("J6").Value = 10.11286928358157
$ws.Range("K6").Value = 14.10316019708641
$ws.Range("N6").Value = 19.60510113163344
$ws.Range("B7").Value = 7.699035920484351
$ws.Range("D7").Value = 8.102359513287801
$ws.Range("E7").Value = 13.17791229222871
$ws.Range("F7").Value = 37.6687024741587
$ws.Range("G7").Value = 42.88685471771691
$ws.Range("H7").Value = 17.72117562294613
$ws.Range("J7").Value = 10.10937933872252
$ws.Range("K7").Value = 14.19877900726222
$ws.Range("N7").Value = 19.58332401942838
$ws.Range("B8").Value = 7.787539156245078
$ws.Range("D8").Value = 8.14218857075919
$ws.Range("E8").Value = 13.23517527128372
$ws.Range("F8").Value = 37.71219198222902
$ws.Range("G8").Value = 43.03153581922303
$ws.Range("H8").Value = 17.67350796742829
$ws.Range("J8").Value = 10.09776785727069
$ws.Range("K8").Value = 14.62010042324182
$ws.Range("N8").Value = 19.49153527237102
$ws.Range("B9").Value = 7.972389438921705
$ws.Range("D9").Value = 8.229715036511925
$ws.Range("E9").Value = 13.36877658529951
$ws.Range("F9").Value = 37.87904964918308
$ws.Range("G9").Value = 43.42648358390441
$ws.Range("H9").Value = 17.60918624373828
$ws.Range("J9").Value = 10.08861529709966
$ws.Range("K9").Value = 15.43613163543179
$ws.Range("N9").Value = 19.32691190180327
$ws.Range("B10").Value = 8.113889286877505
$ws.Range("D10").Value = 8.299417135780441
$ws.Range("E10").Value = 13.47938857571292
$ws.Range("F10").Value = 38.05099102587152
$ws.Range("G10").Value = 43.7831044148108
$ws.Range("H10").Value = 17.57981870583207
$ws.Range("J10").Value = 10.09018522498903
$ws.Range("K10").Value = 16.02282386482859
$ws.Range("N10").Value = 19.21522635172461
$ws.Range("B11").Value = 8.179235102325572
$ws.Range("D11").Value = 8.332218271811339
$ws.Range("E11").Value = 13.53228433472409
$ws.Range("F11").Value = 38.13979827202828
$ws.Range("G11").Value = 43.95938337660627
$ws.Range("H11").Value = 17.57036397490628
$ws.Range("J11").Value = 10.09269633633719
$ws.Range("K11").Value = 16.2856495239298
$ws.Range("N11").Value = 19.16640573491161
$ws.Range("B12").Value = 8.204098171146255
$ws.Range("D12").Value = 8.344789258705417
$ws.Range("E12").Value = 13.55267332941566
$ws.Range("F12").Value = 38.17493533109441
$ws.Range("G12").Value = 44.02811653538745
$ws.Range("H12").Value = 17.56734650115915
$ws.Range("J12").Value = 10.09390501255954
$ws.Range("K12").Value = 16.38449216429364
$ws.Range("N12").Value = 19.14820243424348
$ws.Range("B13").Value = 8.198738571972306
$ws.Range("D13").Value = 8.342075323237715
$ws.Range("E13").Value = 13.54826644946326
$ws.Range("F13").Value = 38.16730114825337
$ws.Range("G13").Value = 44.01322623815515
$ws.Range("H13").Value = 17.56797131748197
$ws.Range("J13").Value = 10.0936332494456
$ws.Range("K13").Value = 16.3632365509555
$ws.Range("N13").Value = 19.15211023375522
$ws.Range("B14").Value = 8.181278364918438
$ws.Range("D14").Value = 8.333249532721748
$ws.Range("E14").Value = 13.53395463615457
$ws.Range("F14").Value = 38.14265889794289
$ws.Range("G14").Value = 43.96499865552145
$ws.Range("H14").Value = 17.5701044388255
$ws.Range("J14").Value = 10.0927906127868
$ws.Range("K14").Value = 16.29379561505328
$ws.Range("N14").Value = 19.16490245523348
$ws.Range("B15").Value = 8.170598189120588
$ws.Range("D15").Value = 8.327862783031323
$ws.Range("E15").Value = 13.5252345556316
$ws.Range("F15").Value = 38.12776066304978
$ws.Range("G15").Value = 43.93571447457939
$ws.Range("H15").Value = 17.57148436955885
$ws.Range("J15").Value = 10.09230802217951
$ws.Range("K15").Value = 16.25116916335034
$ws.Range("N15").Value = 19.17277500103337
$ws.Range("B16").Value = 8.109636521809172
$ws.Range("D16").Value = 8.29729488501798
$ws.Range("E16").Value = 13.47598255939855
$ws.Range("F16").Value = 38.04539912114563
$ws.Range("G16").Value = 43.77186355210974
$ws.Range("H16").Value = 17.58051528703654
$ws.Range("J16").Value = 10.09005722679356
$ws.Range("K16").Value = 16.00555716536428
$ws.Range("N16").Value = 19.21845671987217
$ws.Range("B17").Value = 8.07247273567604
$ws.Range("D17").Value = 8.278817365006402
$ws.Range("E17").Value = 13.44641983220757
$ws.Range("F17").Value = 37.99757533169399
$ws.Range("G17").Value = 43.6749176408724
$ws.Range("H17").Value = 17.58705657753766
$ws.Range("J17").Value = 10.08913627656966
$ws.Range("K17").Value = 15.85376846337588
$ws.Range("N17").Value = 19.24698846582279
$ws.Range("B18").Value = 8.051190398380761
$ws.Range("D18").Value = 8.268292943532629
$ws.Range("E18").Value = 13.42965953768969
$ws.Range("F18").Value = 37.97106621418268
$ws.Range("G18").Value = 43.62048209059417
$ws.Range("H18").Value = 17.59118641940455
$ws.Range("J18").Value = 10.08877575345796
$ws.Range("K18").Value = 15.76608664171767
$ws.Range("N18").Value = 19.26358616594306
$ws.Range("B19").Value = 8.044001271080164
$ws.Range("D19").Value = 8.264747526184216
$ws.Range("E19").Value = 13.42402695555536
$ws.Range("F19").Value = 37.96226249103609
$ws.Range("G19").Value = 43.60227990656495
$ws.Range("H19").Value = 17.59264778078346
$ws.Range("J19").Value = 10.08868276064853
$ws.Range("K19").Value = 15.7363374798088
$ws.Range("N19").Value = 19.26923802542567
$ws.Range("B20").Value = 8.076419391654541
$ws.Range("D20").Value = 8.280773683469295
$ws.Range("E20").Value = 13.44954172977879
$ws.Range("F20").Value = 38.00256308340382
$ws.Range("G20").Value = 43.68510082642609
$ws.Range("H20").Value = 17.58632220431215
$ws.Range("J20").Value = 10.08921680756028
$ws.Range("K20").Value = 15.86996638732061
$ws.Range("N20").Value = 19.24393187064333
$ws.Range("B21").Value = 8.186403826394494
$ws.Range("D21").Value = 8.335837872198438
$ws.Range("E21").Value = 13.53814873248232
$ws.Range("F21").Value = 38.14985613848852
$ws.Range("G21").Value = 43.97911086681921
$ws.Range("H21").Value = 17.56946260494004
$ws.Range("J21").Value = 10.09303112539865
$ws.Range("K21").Value = 16.31421141513882
$ws.Range("N21").Value = 19.16113737399492
$ws.Range("B22").Value = 8.258962608390272
$ws.Range("D22").Value = 8.372695743654363
$ws.Range("E22").Value = 13.59814219848486
$ws.Range("F22").Value = 38.25489994478861
$ws.Range("G22").Value = 44.18278363497204
$ws.Range("H22").Value = 17.56172492177231
$ws.Range("J22").Value = 10.09702611065866
$ws.Range("K22").Value = 16.60052472238505
$ws.Range("N22").Value = 19.10868109143755
$ws.Range("B23").Value = 8.220182017827755
$ws.Range("D23").Value = 8.352946847058368
$ws.Range("E23").Value = 13.56593607229402
$ws.Range("F23").Value = 38.19803838873317
$ws.Range("G23").Value = 44.07303991801319
$ws.Range("H23").Value = 17.56555407460506
$ws.Range("J23").Value = 10.09475671178812
$ws.Range("K23").Value = 16.44811338322625
$ws.Range("N23").Value = 19.13652709525542
$ws.Range("B24").Value = 8.074634846324033
$ws.Range("D24").Value = 8.279888924378996
$ws.Range("E24").Value = 13.44812958477551
$ws.Range("F24").Value = 38.00030505008139
$ws.Range("G24").Value = 43.6804929558165
$ws.Range("H24").Value = 17.58665306453129
$ws.Range("J24").Value = 10.08917987320819
$ws.Range("K24").Value = 15.86264459545089
$ws.Range("N24").Value = 19.24531315155215
$ws.Range("B25").Value = 7.921295072340547
$ws.Range("D25").Value = 8.205062806834363
$ws.Range("E25").Value = 13.33040495088877
$ws.Range("F25").Value = 37.82520686993696
$ws.Range("G25").Value = 43.30785305699841
$ws.Range("H25").Value = 17.62345295597995
$ws.Range("J25").Value = 10.08963368852625
$ws.Range("K25").Value = 15.21717127514022
$ws.Range("N25").Value = 19.36981230883078
